$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(
    731.634050990736,
    607.870083681745,
    614.870684879077,
    823.970611832034,
    916.025500464003,
    1010.62988789579,
    1051.80736031454,
    1007.45624683723,
    874.850218908653,
    779.681175002933,
    629.47090436452,
    578.766158452991,
    749.289668387384,
    736.180348042971,
    657.74419003204,
    773.609509520636,
    857.713982962948,
    960.867376738233,
    1011.14164987985,
    981.392891935124,
    868.096097779587,
    766.308935591763,
    619.276558587264,
    580.909522634798,
    758.969991767529,
    808.709524577327,
    695.46451195807,
    714.127001191491,
    809.940366914215,
    893.474441349911,
    965.927832167815,
    957.727181149081,
    865.744685987109,
    759.824791357511,
    614.088313255905,
    581.716724469889,
    763.691181483323,
    840.372984954311,
    724.900173277365,
    655.000810551545,
    771.077869928791,
    850.360146738297,
    936.265091577122,
    946.085520695474,
    867.682338716199,
    756.234563338431,
    611.779902750961,
    582.364904818331
)

$newI = @(
    218.634050990736,
    243.870083681745,
    254.870684879077,
    408.970611832034,
    309.025500464003,
    178.629887895792,
    221.807360314541,
    177.456246837231,
    26.8502189086531,
    90.6811750029334,
    112.47090436452,
    95.7661584529911,
    387.289668387384,
    409.180348042971,
    170.74419003204,
    251.609509520636,
    208.713982962948,
    131.867376738233,
    72.1416498798529,
    68.3928919351238,
    82.0960977795871,
    111.308935591763,
    95.2765585872639,
    93.9095226347982,
    353.969991767529,
    528.709524577327,
    287.46451195807,
    293.127001191491,
    243.940366914215,
    193.474441349911,
    224.927832167815,
    227.727181149081,
    234.744685987109,
    265.824791357511,
    195.088313255905,
    347.716724469889,
    579.691181483323,
    499.372984954311,
    272.900173277365,
    203.000810551545,
    224.077869928791,
    223.360146738297,
    279.265091577122,
    268.085520695474,
    271.682338716199,
    276.234563338431,
    234.779902750961,
    224.364904818331
)

for ($idx = 0; $idx -lt $newB.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 2).Value = $newB[$idx]
    $ws.Cells.Item($row, 9).Value = $newI[$idx]
}
